# "permission to dance, hard 노트 수정"
#
# The workbook is a rhythm-game note chart. Each data row has columns
# D,E,F,G representing four "lanes"; a 1 (highlighted with fill style
# index 1) marks a note in that lane for that row, 0 (default/no fill)
# means no note. This edit moves a batch of notes around in rows
# 360-808 (toggling specific D/E/F/G cells between "note" and "no note"),
# and updates the current selection in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that go from "no note" (value 0, default style) to
# "note" (value 1, fill style).
$toOne = $ws.Range("E496,E512,D552,E560,F564,F568,D572,E576,F580,F584,E592,D604,E608,G608,F612,G620,F636,F644,D648,D664,F680,E684,F696,F700,G704,F712,E716,G720,G736,F740,E748,G752,D756,F764,G768,D772,F772,F776,E788,G800,F804,D805,E806,F807,G808")

# Cells that go from "note" (value 1, fill style) back to
# "no note" (value 0, default style).
$toZero = $ws.Range("D360,D480,E552,G552,G560,E564,G564,G568,E572,F574,G576,D580,G580,E586,G592,D596,D600,E602,E612,G618,D632,F634,F640,F648,G650,G680,E682,G684,G692,G696,E700,F702,E704,E714,G716,F720,F724,G734,F736,G746,E752,F756,F760,F762,E768,F768,F780,G786,D788,E800,D804,D808")

# D360 still carries the "note" fill style at this point, so use it as
# the format source for every cell that needs to turn into a note.
# (Its own value/style gets reset afterwards via the $toZero pass.)
$template = $ws.Range("D360")
$template.Copy()
foreach ($area in $toOne.Areas) {
    $area.PasteSpecial(-4122)
    $area.Value2 = 1
}

foreach ($area in $toZero.Areas) {
    $area.ClearFormats()
    $area.Value2 = 0
}

# Update the active cell / selection to match the new editing position.
$ws.Range("H806").Select()
